$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.300.57'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '3.497.27'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.68%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.484'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.70'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.62%  '
$ws.Range("E10").Value = '  -0.44%  '
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("D12").Value = '4.092.35'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("E13").Value = '  +0.02%  '
$ws.Range("E14").Value = '  -1.21%  '
$ws.Range("D15").Value = '3.496.60'
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").Value = '64.239.40'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.38'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.49%  '
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("E19").Value = '  -0.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.49'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '385.63'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.578'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.44%  '
$ws.Range("D23").Value = '3.635.98'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("E27").Value = '  +1.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.993'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.98%  '
$ws.Range("E29").Value = '  -2.43%  '
$ws.Range("E30").Value = '  +0.28%  '
$ws.Range("E31").Value = '  +0.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.07'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.154'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.47%  '
$ws.Range("D34").Value = '3.524.79'
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.23'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.38'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.10%  '
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("E39").Value = '  -1.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '164.12'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0781'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("E42").Value = '  +0.20%  '
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("E44").Value = '  -0.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.19'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '24.13'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.72%  '
$ws.Range("E47").Value = '  -0.81%  '
$ws.Range("D48").Value = '2.419.66'
$ws.Range("E48").Value = '  -2.41%  '
$ws.Range("B49").Value = 'SuiNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.923'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.74%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.79'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.13%  '
$ws.Range("E51").Value = '  -1.55%  '
